$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "66.898.32"
Set-TextValue "E2" "  +1.52%  "
Set-TextValue "D3" "3.906.82"
Set-TextValue "E3" "  +3.29%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "469.15"
Set-TextValue "E5" "  +10.15%  "
Set-TextValue "D6" "145.32"
Set-TextValue "E6" "  +6.16%  "
Set-TextValue "E7" "  +2.80%  "
Set-TextValue "E8" "  -0.12%  "
Set-TextValue "D9" "0.742"
Set-TextValue "E9" "  +0.39%  "
Set-TextValue "D10" "0.165"
Set-TextValue "E10" "  +4.28%  "
Set-TextValue "D11" "0.0000338"
Set-TextValue "E11" "  +3.34%  "
Set-TextValue "D12" "43.19"
Set-TextValue "E12" "  -0.09%  "
Set-TextValue "D13" "10.44"
Set-TextValue "E13" "  -2.06%  "
Set-TextValue "D14" "4.528.63"
Set-TextValue "E14" "  +3.24%  "
Set-TextValue "D15" "15.00"
Set-TextValue "E15" "  -0.55%  "
Set-TextValue "D16" "3.890.97"
Set-TextValue "E16" "  +2.15%  "
Set-TextValue "E17" "  -0.37%  "
Set-TextValue "D18" "20.02"
Set-TextValue "E18" "  -0.10%  "
Set-TextValue "E19" "  +3.48%  "
Set-TextValue "D20" "67.204.49"
Set-TextValue "E20" "  +1.69%  "
Set-TextValue "D21" "432.48"
Set-TextValue "E21" "  +5.36%  "
Set-TextValue "D22" "14.74"
Set-TextValue "E22" "  -3.05%  "
Set-TextValue "D23" "3.36"
Set-TextValue "E23" "  +2.69%  "
Set-TextValue "D24" "88.56"
Set-TextValue "E24" "  +3.35%  "
Set-TextValue "D25" "38.74"
Set-TextValue "E25" "  +4.90%  "
Set-TextValue "D26" "3.54"
Set-TextValue "E26" "  +6.98%  "
Set-TextValue "D27" "10.14"
Set-TextValue "E27" "  +3.52%  "
Set-TextValue "E28" "  +4.63%  "
Set-TextValue "E29" "  -1.91%  "
Set-TextValue "D30" "735.52"
Set-TextValue "E30" "  +4.29%  "
Set-TextValue "E31" "  -2.28%  "
Set-TextValue "E32" "  +3.17%  "
Set-TextValue "E33" "  +0.00%  "
Set-TextValue "D34" "43.98"
Set-TextValue "E34" "  +10.10%  "
Set-TextValue "D35" "0.158"
Set-TextValue "E35" "  +3.95%  "
Set-TextValue "D36" "58.10"
Set-TextValue "E36" "  +3.74%  "
Set-TextValue "E37" "  +0.13%  "
Set-TextValue "B38" "ThetaToken"
Set-TextValue "C38" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D38" "3.28"
Set-TextValue "E38" "  +13.76%  "
Set-TextValue "B39" "NEARProtocol"
Set-TextValue "C39" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "5.36"
Set-TextValue "E39" "  -7.37%  "
Set-TextValue "E40" "  +0.61%  "
Set-TextValue "D41" "0.0₃0745"
Set-TextValue "E41" "  +9.23%  "
Set-TextValue "E42" "  +1.01%  "
Set-TextValue "E43" "  +4.47%  "
Set-TextValue "D45" "2.81"
Set-TextValue "E45" "  +6.90%  "
Set-TextValue "E46" "  +6.32%  "
Set-TextValue "E47" "  +1.84%  "
Set-TextValue "D48" "2.48"
Set-TextValue "E48" "  -5.90%  "
Set-TextValue "D49" "3.17"
Set-TextValue "E49" "  -0.17%  "
Set-TextValue "D50" "26.66"
Set-TextValue "E50" "  +3.73%  "
Set-TextValue "D51" "2.88"
Set-TextValue "E51" "  +2.61%  "
